# Update data for regression generation (preprod)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New policy number / claim date values for rows 15-17 (NroPoliza / FechaSiniestro)
# Leading apostrophe forces text entry so the existing text-style/quote-prefix
# formatting of these cells is preserved (values look numeric/date-like).
$ws.Range("F15").Value = "'04104017203 "
$ws.Range("H15").Value = "'14/07/2021"

$ws.Range("F16").Value = "'04104017203 "
$ws.Range("H16").Value = "'14/07/2021"

$ws.Range("F17").Value = "'04104017203 "
$ws.Range("H17").Value = "'14/07/2021"

# Flip the "CLEAS" (Si/No) flag for rows 15-17
$ws.Range("AH15").Value = "No"
$ws.Range("AH16").Value = "No"
$ws.Range("AH17").Value = "Sí"

# Update the active selection on the sheet
$ws.Range("A15").Select()
